$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new value for column F ("想去人数")
$updates = @{
    "展览" = @{
        2  = 592
        4  = 1271
        5  = 1120
        6  = 14145
        7  = 15794
        11 = 194
        20 = 1224
        23 = 16
        24 = 6238
        27 = 5590
        30 = 128
        31 = 4578
    }
    "全部类型" = @{
        2  = 592
        4  = 1271
        5  = 1120
        6  = 14145
        7  = 15794
        11 = 194
        20 = 1224
        24 = 16
        25 = 6238
        28 = 5590
        31 = 128
        32 = 4578
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
